$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '65.162.25'
$ws.Cells.Item(2, 5).Value = '  -5.43%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.387.90'
$ws.Cells.Item(3, 5).Value = '  -6.92%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.20%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '186.88'
$ws.Cells.Item(5, 5).Value = '  -8.54%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '529.40'
$ws.Cells.Item(6, 5).Value = '  -7.56%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.602'
$ws.Cells.Item(7, 5).Value = '  -4.08%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.385.64'
$ws.Cells.Item(8, 5).Value = '  -6.86%  '
$ws.Cells.Item(9, 5).Value = '  +0.05%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.626'
$ws.Cells.Item(10, 5).Value = '  -8.71%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '58.35'
$ws.Cells.Item(11, 5).Value = '  -6.25%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.133'
$ws.Cells.Item(12, 5).Value = '  -12.96%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000254'
$ws.Cells.Item(13, 5).Value = '  -12.96%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '9.29'
$ws.Cells.Item(14, 5).Value = '  -8.56%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.912.32'
$ws.Cells.Item(15, 5).Value = '  -7.24%  '
$ws.Cells.Item(16, 5).Value = '  -3.16%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.381.21'
$ws.Cells.Item(17, 5).Value = '  -6.86%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '64.864.46'
$ws.Cells.Item(18, 5).Value = '  -5.52%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.42'
$ws.Cells.Item(19, 5).Value = '  -8.93%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '11.15'
$ws.Cells.Item(20, 5).Value = '  -10.83%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.972'
$ws.Cells.Item(21, 5).Value = '  -10.64%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '372.62'
$ws.Cells.Item(22, 5).Value = '  -8.63%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '81.23'
$ws.Cells.Item(23, 5).Value = '  -5.69%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.73'
$ws.Cells.Item(24, 5).Value = '  -11.91%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '10.87'
$ws.Cells.Item(25, 5).Value = '  -17.88%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.77'
$ws.Cells.Item(26, 5).Value = '  -5.02%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '5.82'
$ws.Cells.Item(27, 5).Value = '  -5.55%  '
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '11.51'
$ws.Cells.Item(28, 5).Value = '  -9.73%  '
$ws.Cells.Item(29, 2).Value = 'ImmutableX'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.65'
$ws.Cells.Item(29, 5).Value = '  -10.94%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '8.54'
$ws.Cells.Item(30, 5).Value = '  -9.63%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '29.62'
$ws.Cells.Item(31, 5).Value = '  -7.11%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '667.74'
$ws.Cells.Item(32, 5).Value = '  -0.36%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.73'
$ws.Cells.Item(33, 5).Value = '  -17.60%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '11.21'
$ws.Cells.Item(34, 5).Value = '  -9.67%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '61.08'
$ws.Cells.Item(35, 5).Value = '  -4.53%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.106'
$ws.Cells.Item(36, 5).Value = '  -9.10%  '
$ws.Cells.Item(37, 5).Value = '  -0.13%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '36.68'
$ws.Cells.Item(38, 5).Value = '  -13.95%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.381'
$ws.Cells.Item(39, 5).Value = '  -10.43%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.996'
$ws.Cells.Item(40, 5).Value = '  -0.16%  '
$ws.Cells.Item(41, 5).Value = '  -7.33%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.851.81'
$ws.Cells.Item(42, 5).Value = '  -11.67%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.75'
$ws.Cells.Item(43, 5).Value = '  -15.69%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.64'
$ws.Cells.Item(44, 5).Value = '  -8.52%  '
$ws.Cells.Item(45, 2).Value = 'PEPE'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0₃0627'
$ws.Cells.Item(45, 5).Value = '  -22.80%  '
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0393'
$ws.Cells.Item(46, 5).Value = '  -7.08%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.34'
$ws.Cells.Item(47, 5).Value = '  -15.58%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '138.13'
$ws.Cells.Item(48, 5).Value = '  -0.88%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.125'
$ws.Cells.Item(49, 5).Value = '  -5.88%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.87'
$ws.Cells.Item(50, 5).Value = '  -6.89%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.59'
$ws.Cells.Item(51, 5).Value = '  -6.79%  '
